# Generate Report for Handoff
# Updates the localization status report:
#  - Status text "Handed back: in sync with en-US" -> "Ready for handoff"
#  - Refreshed generation timestamps
#  - Narrower "Status" columns to fit the shorter text

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-15 18:53:53"

$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-15 18:53:48"

$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797

# --- de-de sheet ------------------------------------------------------
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-15 18:53:53"

$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
